$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GermanyAllNumbers")

# ---------------------------------------------------------------------------
# 1. Append 27 new daily rows (450-476) of German COVID numbers, continuing
#    the existing pattern used by the preceding rows.
# ---------------------------------------------------------------------------

# Reuse the number formatting (date format for column D, the two decimal
# formats for columns G and J) already present on row 449 so no new cell
# styles are introduced.
$ws.Range("D449:J449").Copy()
$ws.Range("D450:J476").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Repeated iso_code / continent / location text columns.
$ws.Range("A450:A476").Value = "DEU"
$ws.Range("B450:B476").Value = "Europe"
$ws.Range("C450:C476").Value = "Germany"

# Columns K through P stay at 0 for every one of the new rows.
$ws.Range("K450:P476").Value = 0

# Day-by-day new_tests (F) and total_tests (H) numbers copied from the
# published update.
$newTests = @{
    450 = 2440; 451 = 1117; 452 = 1204; 453 = 3254; 454 = 3187; 455 = 2440
    456 = 1911; 457 = 1489; 458 = 549;  459 = 652;  460 = 1455; 461 = 1330
    462 = 1076;             464 = 842;  465 = 346;  466 = 455;  467 = 1016
    468 = 1008; 469 = 774;  470 = 592;  471 = 538;  472 = 219;  473 = 404
    474 = 808;  475 = 892;  476 = 649
}
$totalTests = @{
    450 = 89222; 451 = 89244; 452 = 89384; 453 = 89491; 454 = 89585
    455 = 89687; 456 = 89816; 457 = 89834; 458 = 89844; 459 = 89937
    460 = 90074; 461 = 90179; 462 = 90270; 463 = 90270; 464 = 90385
    465 = 90395; 466 = 90472; 467 = 90523; 468 = 90616; 469 = 90678
    470 = 90746; 471 = 90754; 472 = 90762; 473 = 90819; 474 = 90875
    475 = 90938; 476 = 91007
}

foreach ($r in 450..476) {
    $prev = $r - 1

    # date (column D) - shared "+1 day" pattern
    $ws.Range("D$r").Formula = "=D$prev+1"

    # new_tests (column F) - literal, except row 463 which is derived
    if ($r -eq 463) {
        $ws.Range("F463").Formula = "=E463-E462"
    } else {
        $ws.Range("F$r").Value = $newTests[$r]
    }

    # total_cases (column E) - cumulative running total, except the two
    # rows around the mid-week correction (463/464) which differ.
    if ($r -eq 463) {
        $ws.Range("E463").Formula = "=E464-F464"
    } elseif ($r -eq 464) {
        $ws.Range("E464").Value = 3721981
    } else {
        $ws.Range("E$r").Formula = "=E$prev+F$r"
    }

    # 7-day average of new_tests (column G)
    $g1 = $r - 6
    $ws.Range("G$r").Formula = "=SUM(F$g1`:F$r)/7"

    # total_tests (column H) - literal
    $ws.Range("H$r").Value = $totalTests[$r]

    # new daily test delta (column I)
    $ws.Range("I$r").Formula = "=H$r-H$prev"

    # 7-day average of the daily delta (column J)
    $j1 = $r - 6
    $ws.Range("J$r").Formula = "=SUM(I$j1`:I$r)/7"
}

# ---------------------------------------------------------------------------
# 2. Restore the view state Excel persists after scrolling to the new tail
#    of the sheet (dimension grows automatically with the data above).
# ---------------------------------------------------------------------------
$ws.Range("H477").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 446
